# IATI partner org extract simplified
# Adds a new row of data (row 23) to the "IATI activity IDs" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IATI activity IDs")

$ws.Range("A23").Value = "GB-COH-877338-GV-GOV-1-300708-124"
$ws.Range("B23").Value = "Foreign, Commonwealth and Development Office"
$ws.Range("C23").Value = "GB-GOV-1-300708-124"
$ws.Range("D23").Value = "FCDO Research - Programmes"
$ws.Range("E23").Value = "Institute of Development Studies"
